# Updated cryptos list — applies the latest price/volume snapshot to
# Sheet1 columns D (Price) and E (Volume(1h)) for rows 2-51.
#
# Source values are plain text (not numbers) in the original workbook —
# e.g. "63.412.42" / "549.38" — some of which are lexically numeric-looking
# ("549.38", "1.00", ...). Assigning such a string straight to .Value lets
# Excel's COM layer auto-coerce it to a real Number cell, which would both
# change the stored type and perturb its formatting (e.g. "1.00" -> 1).
# To keep every written cell a literal text cell (matching the source
# file's inlineStr cells) we momentarily force a Text number format before
# the write, then call ClearFormats() to drop the now-unneeded explicit
# format back to the sheet default — leaving only the cell's text content
# changed, exactly as in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.379.50'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.80%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.065.22'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.22%  '
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.40'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.15'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.14%  '
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.059.73'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("E8").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.55%  '
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.50'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.77%  '
$ws.Range("E10").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("E11").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000227'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.96%  '
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.83'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.16%  '
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.558.52'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.262.20'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.061.32'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("E17").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("E18").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.27%  '
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '481.85'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.48%  '
$ws.Range("E20").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.673'
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.23'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.13%  '
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.69'
$ws.Range("D24").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.60'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.85%  '
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E26").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("E27").ClearFormats()

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.92'
$ws.Range("D28").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E28").ClearFormats()

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.09%  '
$ws.Range("E29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("E30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.09'
$ws.Range("D31").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.43%  '
$ws.Range("E31").ClearFormats()

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.45'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +7.23%  '
$ws.Range("E33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.71'
$ws.Range("D34").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.67%  '
$ws.Range("E34").ClearFormats()

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.52'
$ws.Range("D35").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E35").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.89%  '
$ws.Range("E36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '469.56'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.27%  '
$ws.Range("E37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0819'
$ws.Range("D38").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0397'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.19%  '
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.072.15'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.87%  '
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.119'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("E41").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.58'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.52%  '
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("E44").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.54%  '
$ws.Range("E45").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.21%  '
$ws.Range("E47").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.61'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.15%  '
$ws.Range("E49").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.15%  '
$ws.Range("E50").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.26%  '
$ws.Range("E51").ClearFormats()
